$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: table name & description ---
$ws.Range("B1").Value = "Vagas"
$ws.Range("B2").Value = "Tabela responsável por armazenar os dados das vagas"

# --- Attribute rows 5-6: simple renames, keep existing formatting ---
$ws.Range("A5").Value = "id"
$ws.Range("H5").Value = "Código de identificador das vagas"

$ws.Range("A6").Value = "numero"
$ws.Range("H6").Value = "Número denominado a vaga"

# --- Row 7 (was "andar"/tinyint) becomes "tipo_id" FK column ---
# Bring over the distinctive left/right-only border used on C8 so the
# new FK column matches the rest of the table visually.
$ws.Range("C8").Copy($ws.Range("C7"))
$ws.Range("A7").Value = "tipo_id"
$ws.Range("C7").Value = "int"
$ws.Range("D7").Value = "1 – sem limite"
$ws.Range("E7").ClearContents()
$ws.Range("G7").Value = "X"
$ws.Range("H7").Value = "Foreign Key da tabela TipoVaga"

# --- Row 8 (was "tipo"/varchar) becomes "veiculo_id" FK column ---
$ws.Range("A8").Value = "veiculo_id"
$ws.Range("C8").Value = "int"
$ws.Range("D8").Value = "1 – sem limite"
$ws.Range("E8").ClearContents()
$ws.Range("G8").Value = "X"
$ws.Range("H8").Value = "Foreign Key da tabela Veículos"

# --- Row 9 was a blank spacer row; turn it into the "disponivel" column,
#     copying formatting from the row above (row 6) first. ---
$ws.Range("A6:H6").Copy($ws.Range("A9:H9"))
$ws.Range("A9").Value = "disponivel"
$ws.Range("C9").Value = "tinyint"
$ws.Range("D9").Value = "0 – 1"
$ws.Range("E9").Value = "NOT NULL"
$ws.Range("F9").ClearContents()
$ws.Range("H9").Value = "Identificação da disponibilidade da vaga"

# --- Index table rows 13-15 ---
$ws.Range("A13").Value = "PRIMARY"
$ws.Range("C13").Value = "Sim"
$ws.Range("D13").Value = "Não"
$ws.Range("E13").Value = "Sim"
$ws.Range("F13").Value = "id"

$ws.Range("A14").Value = "INDEX_tipo_id"
$ws.Range("C14").Value = "Não"
$ws.Range("D14").Value = "Sim"
$ws.Range("E14").Value = "Não"
$ws.Range("F14").Value = "tipo_id"

# Row 15 previously carried a stray underlined/gray style on F:H; align its
# formatting with row 14 before writing the new values.
$ws.Range("A14:H14").Copy($ws.Range("A15:H15"))
$ws.Range("A15").Value = "INDEX_veiculo_id"
$ws.Range("C15").Value = "Não"
$ws.Range("D15").Value = "Sim"
$ws.Range("E15").Value = "Não"
$ws.Range("F15").Value = "veiculo_id"

# --- Column width & selection, matching the final authoring state ---
$ws.Columns.Item(8).ColumnWidth = 68.88671875
$ws.Range("F15:H15").Select()
